# Bugfixed battery and wifi level on Cloud API
# Adds the two new ToDo rows that track the Cloud API work:
#   - "Add percentage widgets for battery and wifi level in app" -> Closed
#   - "Sync app to latest SDK"                                    -> Open
# and moves the sheet's view/selection down to the newly added rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 63: battery/wifi percentage widgets, status Closed
$ws.Range("B63").Value = "Add percentage widgets for battery and wifi level in app"
$ws.Range("C63").Value = "Closed"

# New row 64: sync app to latest SDK, status Open
$ws.Range("B64").Value = "Sync app to latest SDK"
$ws.Range("C64").Value = "Open"

# Scroll the view towards the new rows and select the last added status cell,
# matching where the author was working when they saved the file.
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("C63").Select()
